# Generate Report for Handoff
# Update the handoff UUID/file names and timestamps across the Overview,
# zh-cn and de-de sheets (both the cell values and the matching hyperlink
# display text).

$wb = $excel.ActiveWorkbook

$oldId = "3b02391b-a675-468b-b827-8e5ab77c5869"
$newId = "ba81ee4b-15d8-4cf6-8ce8-ee74b9ee42c5"

$oldMd = "$oldId.md"
$newMd = "$newId.md"

$oldZh = "$oldId.1f204e2fdcc5ee6a3b42b0dcf2128e76a7950c44.zh-cn.xlf"
$newZh = "$newId.e35368e0e97e1f37c512f5a19e2900caf56b94f6.zh-cn.xlf"

$oldDe = "$oldId.1f204e2fdcc5ee6a3b42b0dcf2128e76a7950c44.de-de.xlf"
$newDe = "$newId.e35368e0e97e1f37c512f5a19e2900caf56b94f6.de-de.xlf"

$oldOverviewDate = "2016-09-13 23:09:18"
$newOverviewDate = "2016-10-13 23:10:19"

$oldZhDate = "2016-03-13 23:06:42"
$newZhDate = "2016-03-13 23:10:15"

$oldDeDate = "2016-03-13 23:09:18"
$newDeDate = "2016-03-13 23:10:19"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("D2").Value = $newOverviewDate

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMd) {
        $hl.TextToDisplay = $newMd
    }
}

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("D2").Value = $newZh
$wsZh.Range("E2").Value = $newZhDate

foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMd) {
        $hl.TextToDisplay = $newMd
    } elseif ($hl.TextToDisplay -eq $oldZh) {
        $hl.TextToDisplay = $newZh
    }
}

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("D2").Value = $newDe
$wsDe.Range("E2").Value = $newDeDate

foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMd) {
        $hl.TextToDisplay = $newMd
    } elseif ($hl.TextToDisplay -eq $oldDe) {
        $hl.TextToDisplay = $newDe
    }
}
